$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header cells: A1 "Zip" -> "zip_code", B1 "Encoding" -> "big_city"
$ws.Range("A1").Value = "zip_code"
$ws.Range("B1").Value = "big_city"

# Column C ("column" header + "big_city" filler values in every data row) is
# wiped out entirely, leaving empty (but still styled) cells C1:C324.
$ws.Range("C1:C324").ClearContents()

# Move the active selection to C4 (matches the saved cursor position).
$ws.Range("C4").Select() | Out-Null
